$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1728.0938
$ws.Range("I2").Value = 1503.5834
$ws.Range("J2").Value = 2401.625
$ws.Range("K2").Value = 1503.5834
$ws.Range("L2").Value = 2401.625
$ws.Range("M2").Value = -1390.5834
$ws.Range("N2").Value = -2627.625

$ws.Range("H32").Value = 7429.55
$ws.Range("I32").Value = 4128.878
$ws.Range("J32").Value = 22465.945
$ws.Range("K32").Value = 4128.878
$ws.Range("L32").Value = 22465.945
$ws.Range("M32").Value = -3841.878
$ws.Range("N32").Value = -23039.945

$ws.Range("H45").Value = 1196.5
$ws.Range("I45").Value = 706.44446
$ws.Range("K45").Value = 706.44446
$ws.Range("M45").Value = -329.44446

$ws.Range("H61").Value = 34554244
$ws.Range("I61").Value = 45502104
$ws.Range("J61").Value = 146685
$ws.Range("K61").Value = 45502104
$ws.Range("L61").Value = 146685
$ws.Range("M61").Value = -45501892
$ws.Range("N61").Value = -147109

$ws.Range("H74").Value = 7799976
$ws.Range("I74").Value = 9553414
$ws.Range("K74").Value = 9553414
$ws.Range("M74").Value = -9552540

$ws.Range("H77").Value = 7799976
$ws.Range("I77").Value = 9553414
$ws.Range("K77").Value = 47767070
$ws.Range("M77").Value = -47762702

$ws.Range("H110").Value = 1006.0714
$ws.Range("I110").Value = 841.6
$ws.Range("J110").Value = 1417.25
$ws.Range("K110").Value = 841.6
$ws.Range("L110").Value = 1417.25
$ws.Range("M110").Value = 1203.4
$ws.Range("N110").Value = -5507.25

$ws.Range("H116").Value = 1728.0938
$ws.Range("I116").Value = 1503.5834
$ws.Range("J116").Value = 2401.625
$ws.Range("K116").Value = 1503.5834
$ws.Range("L116").Value = 2401.625
$ws.Range("M116").Value = 790.4166
$ws.Range("N116").Value = -6989.625

$ws.Range("H136").Value = 34554244
$ws.Range("I136").Value = 45502104
$ws.Range("J136").Value = 146685
$ws.Range("K136").Value = 136506312
$ws.Range("L136").Value = 440055
$ws.Range("M136").Value = -136503762
$ws.Range("N136").Value = -445155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1728.0938
$ws.Range("I3").Value = 1503.5834
$ws.Range("J3").Value = 2401.625
$ws.Range("K3").Value = 1503.5834
$ws.Range("L3").Value = 2401.625
$ws.Range("M3").Value = -1389.5834
$ws.Range("N3").Value = -2629.625

$ws.Range("H102").Value = 12627.75
$ws.Range("I102").Value = 12627.75
$ws.Range("K102").Value = 12627.75
$ws.Range("M102").Value = -9382.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35729.8
$ws.Range("I31").Value = 18415.062
$ws.Range("J31").Value = 66511.55499999999
$ws.Range("K31").Value = 18415.062
$ws.Range("L31").Value = 66511.55499999999
$ws.Range("M31").Value = -18120.062
$ws.Range("N31").Value = -67101.55499999999

$ws.Range("H34").Value = 35729.8
$ws.Range("I34").Value = 18415.062
$ws.Range("J34").Value = 66511.55499999999
$ws.Range("K34").Value = 18415.062
$ws.Range("L34").Value = 66511.55499999999
$ws.Range("M34").Value = -18213.062
$ws.Range("N34").Value = -66915.55499999999

$ws.Range("H58").Value = 47620570
$ws.Range("I58").Value = 62501548
$ws.Range("J58").Value = 1440
$ws.Range("K58").Value = 62501548
$ws.Range("L58").Value = 1440
$ws.Range("M58").Value = -62501345
$ws.Range("N58").Value = -1846

$ws.Range("H136").Value = 47620570
$ws.Range("I136").Value = 62501548
$ws.Range("J136").Value = 1440
$ws.Range("K136").Value = 187504644
$ws.Range("L136").Value = 4320
$ws.Range("M136").Value = -187502094
$ws.Range("N136").Value = -9420

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 2111.4443
$ws.Range("I47").Value = 1000.75
$ws.Range("K47").Value = 3002.25
$ws.Range("M47").Value = -2571.25

$ws.Range("H68").Value = 668.2632
$ws.Range("I68").Value = 675.0323
$ws.Range("J68").Value = 638.2857
$ws.Range("K68").Value = 2025.0969
$ws.Range("L68").Value = 1914.8571
$ws.Range("M68").Value = -1214.0969
$ws.Range("N68").Value = -3536.8571

$ws.Range("H71").Value = 668.2632
$ws.Range("I71").Value = 675.0323
$ws.Range("J71").Value = 638.2857
$ws.Range("K71").Value = 6075.2907
$ws.Range("L71").Value = 5744.571300000001
$ws.Range("M71").Value = -2019.2907
$ws.Range("N71").Value = -13856.5713

$ws.Range("H136").Value = 3424.3333
$ws.Range("I136").Value = 3111.8
$ws.Range("J136").Value = 3647.5715
$ws.Range("K136").Value = 9335.400000000001
$ws.Range("L136").Value = 10942.7145
$ws.Range("M136").Value = -4235.400000000001
$ws.Range("N136").Value = -21142.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 8020
$ws.Range("I52").Value = 8020
$ws.Range("K52").Value = 8020
$ws.Range("M52").Value = -7761

$ws.Range("H55").Value = 5000
$ws.Range("I55").Value = 5000
$ws.Range("K55").Value = 5000
$ws.Range("M55").Value = -4673

$ws.Range("H122").Value = 1535.3334
$ws.Range("I122").Value = 1476.8334
$ws.Range("J122").Value = 1769.3334
$ws.Range("K122").Value = 4430.5002
$ws.Range("L122").Value = 5308.0002
$ws.Range("M122").Value = -1980.5002
$ws.Range("N122").Value = -10208.0002

$ws.Range("H126").Value = 1346.5454
$ws.Range("I126").Value = 1045.7778
$ws.Range("K126").Value = 3137.3334
$ws.Range("M126").Value = -667.3334000000004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1981.174
$ws.Range("I61").Value = 1713.7894
$ws.Range("K61").Value = 1713.7894
$ws.Range("M61").Value = -1511.7894

$ws.Range("H113").Value = 1981.174
$ws.Range("I113").Value = 1713.7894
$ws.Range("K113").Value = 1713.7894
$ws.Range("M113").Value = 456.2106000000001

$ws.Range("H122").Value = 3296.98
$ws.Range("I122").Value = 2751.1875
$ws.Range("J122").Value = 3553.8235
$ws.Range("K122").Value = 8253.5625
$ws.Range("L122").Value = 10661.4705
$ws.Range("M122").Value = -5803.5625
$ws.Range("N122").Value = -15561.4705

$ws.Range("H136").Value = 36940.76
$ws.Range("I136").Value = 21990.48
$ws.Range("K136").Value = 65971.44
$ws.Range("M136").Value = -63421.44

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 8500
$ws.Range("I58").Value = 8500
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 8500
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -8192
$ws.Range("N58").Value = $null
